# Regenerate the handoff report: refresh the "Latest Handoff Datetime" column
# (D) for every row that is part of the batch that was just handed off, on
# both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhRows = @(7, 10, 11, 12, 13, 14, 15, 16)
$deRows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $wsZh.Cells.Item($r, 4).Value = "2016-03-08 08:19:32"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $wsDe.Cells.Item($r, 4).Value = "2016-03-08 08:19:36"
}
